$wb = $excel.ActiveWorkbook

# --- zh-cn sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Status column -> handed back
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("B3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime (column G) now has real timestamps
$ws.Range("G2").Value = "2016-02-22 17:36:42"
$ws.Range("G3").Value = "2016-02-22 17:36:42"

# Rebuild all hyperlinks (existing ones are recreated, new ones for the
# newly populated "Latest Target File" / "Latest Handback File" columns
# are added) in row-major, left-to-right order.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/328bd31b-a90f-45cb-8f2f-e87ea939f61c.md", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5733f043465cafeb960a231308d229918c4c277/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.zh-cn.xlf", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/328bd31b-a90f-45cb-8f2f-e87ea939f61c.md", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.md")
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5733f043465cafeb960a231308d229918c4c277/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.zh-cn.xlf", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/c67fc2fe-241c-48af-945b-36bb1c92c5b5.md", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5733f043465cafeb960a231308d229918c4c277/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.zh-cn.xlf", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/c67fc2fe-241c-48af-945b-36bb1c92c5b5.md", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.md")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5733f043465cafeb960a231308d229918c4c277/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.zh-cn.xlf", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/.localization-config", "", "", ".localization-config")

# --- de-de sheet -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("de-de")

$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"

$ws2.Range("G2").Value = "2016-02-22 17:37:01"
$ws2.Range("G3").Value = "2016-02-22 17:37:01"

$ws2.Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/328bd31b-a90f-45cb-8f2f-e87ea939f61c.md", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9fc6be741da3967efcd10cbd65c4efc01dde1885/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.de-de.xlf", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.de-de.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/328bd31b-a90f-45cb-8f2f-e87ea939f61c.md", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9fc6be741da3967efcd10cbd65c4efc01dde1885/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.de-de.xlf", "", "", "328bd31b-a90f-45cb-8f2f-e87ea939f61c.e3878a39166d33f5054544308340ffa15ce61a1a.de-de.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/c67fc2fe-241c-48af-945b-36bb1c92c5b5.md", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9fc6be741da3967efcd10cbd65c4efc01dde1885/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.de-de.xlf", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.de-de.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/e2e/c67fc2fe-241c-48af-945b-36bb1c92c5b5.md", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9fc6be741da3967efcd10cbd65c4efc01dde1885/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.de-de.xlf", "", "", "c67fc2fe-241c-48af-945b-36bb1c92c5b5.73cac44c7a2c2e44a811fe66c2536c0cb709451e.de-de.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/54094527cc9aeb58673ec450cd2a0c3d48047d8b/.localization-config", "", "", ".localization-config")
